$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header labels in row 1 (L1, M1)
$ws.Range("L1").Value2 = "add1"
$ws.Range("M1").Value2 = "add2"

# Add the new "G1 LSS" / "G1 Science" columns for the data rows (2-48 only)
$ws.Range("L2:L48").Value2 = "G1 LSS"
$ws.Range("M2:M48").Value2 = "G1 Science"

# Rename column J values from "Normal Technical" to "Lower Secondary NT"
# for every data row (2-50).
$ws.Range("J2:J50").Value2 = "Lower Secondary NT"

# Widen column J slightly to fit the new, longer text
$ws.Columns.Item(10).ColumnWidth = 17.1

# Update the view: move the active selection to C6 (the workbook was also
# scrolled so column C is left-most, but that pane-scroll state isn't
# reachable through the COM surface exposed here)
$ws.Range("C6").Select() | Out-Null
